$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value2 = 11236.111
$ws.Range("I15").Value2 = 11236.111
$ws.Range("K15").Value2 = 33708.333
$ws.Range("M15").Value2 = -33539.333

$ws.Range("H106").Value2 = 60608028
$ws.Range("I106").Value2 = 33335114
$ws.Range("K106").Value2 = 33335114
$ws.Range("M106").Value2 = -33334483

$ws.Range("H112").Value2 = 55574640
$ws.Range("J112").Value2 = 166722670
$ws.Range("L112").Value2 = 500168010
$ws.Range("N112").Value2 = -500170226

$ws.Range("H132").Value2 = 7247626.5
$ws.Range("I132").Value2 = 863.09375
$ws.Range("J132").Value2 = 23811658
$ws.Range("K132").Value2 = 2589.28125
$ws.Range("L132").Value2 = 71434974
$ws.Range("M132").Value2 = -59.28125
$ws.Range("N132").Value2 = -71440034

$ws.Range("H135").Value2 = 2013.75
$ws.Range("I135").Value2 = 1928.2069
$ws.Range("K135").Value2 = 17353.8621
$ws.Range("M135").Value2 = -14818.8621

$ws.Range("H137").Value2 = 1363.775
$ws.Range("I137").Value2 = 1306.4073
$ws.Range("J137").Value2 = 1482.9231
$ws.Range("K137").Value2 = 3919.2219
$ws.Range("L137").Value2 = 4448.7693
$ws.Range("M137").Value2 = -1369.2219
$ws.Range("N137").Value2 = -9548.7693

$ws.Range("H138").Value2 = 3992.6233
$ws.Range("I138").Value2 = 894.2
$ws.Range("J138").Value2 = 6574.643
$ws.Range("K138").Value2 = 2682.6
$ws.Range("L138").Value2 = 19723.929
$ws.Range("M138").Value2 = 2457.4
$ws.Range("N138").Value2 = -30003.929

$ws.Range("H141").Value2 = 2063.7878
$ws.Range("I141").Value2 = 1964.1904
$ws.Range("J141").Value2 = 2238.0833
$ws.Range("K141").Value2 = 5892.5712
$ws.Range("L141").Value2 = 6714.249899999999
$ws.Range("M141").Value2 = -712.5712000000003
$ws.Range("N141").Value2 = -17074.2499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 5182.905
$ws.Range("I32").Value2 = 3599.5134
$ws.Range("K32").Value2 = 3599.5134
$ws.Range("M32").Value2 = -3312.5134

$ws.Range("H61").Value2 = 3722.7173
$ws.Range("I61").Value2 = 3817.1462
$ws.Range("J61").Value2 = 2948.4
$ws.Range("K61").Value2 = 3817.1462
$ws.Range("L61").Value2 = 2948.4
$ws.Range("M61").Value2 = -3605.1462
$ws.Range("N61").Value2 = -3372.4

$ws.Range("H63").Value2 = 111113580
$ws.Range("I63").Value2 = 111113580
$ws.Range("K63").Value2 = 111113580
$ws.Range("M63").Value2 = -111112894

$ws.Range("H66").Value2 = 111113580
$ws.Range("I66").Value2 = 111113580
$ws.Range("K66").Value2 = 555567900
$ws.Range("M66").Value2 = -555564468

$ws.Range("H74").Value2 = 12196297
$ws.Range("I74").Value2 = 1005.71875
$ws.Range("J74").Value2 = 55557336
$ws.Range("K74").Value2 = 1005.71875
$ws.Range("L74").Value2 = 55557336
$ws.Range("M74").Value2 = -131.71875
$ws.Range("N74").Value2 = -55559084

$ws.Range("H77").Value2 = 12196297
$ws.Range("I77").Value2 = 1005.71875
$ws.Range("J77").Value2 = 55557336
$ws.Range("K77").Value2 = 5028.59375
$ws.Range("L77").Value2 = 277786680
$ws.Range("M77").Value2 = -660.59375
$ws.Range("N77").Value2 = -277795416

$ws.Range("H95").Value2 = 39400
$ws.Range("J95").Value2 = 39400
$ws.Range("L95").Value2 = 39400
$ws.Range("N95").Value2 = -44892

$ws.Range("H110").Value2 = 733.3333
$ws.Range("I110").Value2 = 800
$ws.Range("J110").Value2 = 400
$ws.Range("K110").Value2 = 800
$ws.Range("L110").Value2 = 400
$ws.Range("M110").Value2 = 1245
$ws.Range("N110").Value2 = -4490

$ws.Range("H136").Value2 = 3722.7173
$ws.Range("I136").Value2 = 3817.1462
$ws.Range("J136").Value2 = 2948.4
$ws.Range("K136").Value2 = 11451.4386
$ws.Range("L136").Value2 = 8845.200000000001
$ws.Range("M136").Value2 = -8901.438600000001
$ws.Range("N136").Value2 = -13945.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value2 = 3540.8245
$ws.Range("I134").Value2 = 4060.35
$ws.Range("J134").Value2 = 2318.4119
$ws.Range("K134").Value2 = 12181.05
$ws.Range("L134").Value2 = 6955.2357
$ws.Range("M134").Value2 = -9646.049999999999
$ws.Range("N134").Value2 = -12025.2357

$ws.Range("H138").Value2 = 45706
$ws.Range("J138").Value2 = 45706
$ws.Range("L138").Value2 = 45706
$ws.Range("N138").Value2 = -55986

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value2 = 293.79544
$ws.Range("I107").Value2 = 310.14285
$ws.Range("J107").Value2 = 286.16666
$ws.Range("K107").Value2 = 310.14285
$ws.Range("L107").Value2 = 286.16666
$ws.Range("M107").Value2 = 1609.85715
$ws.Range("N107").Value2 = -4126.16666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value2 = 2200.2
$ws.Range("I22").Value2 = 750
$ws.Range("J22").Value2 = 2562.75
$ws.Range("K22").Value2 = 2250
$ws.Range("L22").Value2 = 7688.25
$ws.Range("M22").Value2 = -2081
$ws.Range("N22").Value2 = -8026.25

$ws.Range("H27").Value2 = 2200.2
$ws.Range("I27").Value2 = 750
$ws.Range("J27").Value2 = 2562.75
$ws.Range("K27").Value2 = 2250
$ws.Range("L27").Value2 = 7688.25
$ws.Range("M27").Value2 = -2148
$ws.Range("N27").Value2 = -7892.25

$ws.Range("H44").Value2 = 844.7931
$ws.Range("I44").Value2 = 469.23077
$ws.Range("K44").Value2 = 1407.69231
$ws.Range("M44").Value2 = -1009.69231

$ws.Range("H59").Value2 = 2998
$ws.Range("I59").Value2 = 500
$ws.Range("J59").Value2 = 3622.5
$ws.Range("K59").Value2 = 1500
$ws.Range("L59").Value2 = 10867.5
$ws.Range("M59").Value2 = -960
$ws.Range("N59").Value2 = -11947.5

$ws.Range("H60").Value2 = 378.33334
$ws.Range("I60").Value2 = 254
$ws.Range("K60").Value2 = 762
$ws.Range("M60").Value2 = -511

$ws.Range("H70").Value2 = 1439.579
$ws.Range("I70").Value2 = 489.42856
$ws.Range("J70").Value2 = 4100
$ws.Range("K70").Value2 = 1468.28568
$ws.Range("L70").Value2 = 12300
$ws.Range("M70").Value2 = -1153.28568
$ws.Range("N70").Value2 = -12930

$ws.Range("H73").Value2 = 1439.579
$ws.Range("I73").Value2 = 489.42856
$ws.Range("J73").Value2 = 4100
$ws.Range("K73").Value2 = 1468.28568
$ws.Range("L73").Value2 = 12300
$ws.Range("M73").Value2 = -376.28568
$ws.Range("N73").Value2 = -14484

$ws.Range("H80").Value2 = 3250
$ws.Range("I80").Value2 = 0
$ws.Range("J80").Value2 = 3250
$ws.Range("K80").Value2 = 0
$ws.Range("L80").Value2 = 9750
$ws.Range("N80").Value2 = -11622
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value2 = 3250
$ws.Range("I83").Value2 = 0
$ws.Range("J83").Value2 = 3250
$ws.Range("K83").Value2 = 0
$ws.Range("L83").Value2 = 29250
$ws.Range("N83").Value2 = -38610
$ws.Range("M83").ClearContents()

$ws.Range("H97").Value2 = 11111737
$ws.Range("I97").Value2 = 25000350
$ws.Range("J97").Value2 = 845.2
$ws.Range("K97").Value2 = 75001050
$ws.Range("L97").Value2 = 2535.6
$ws.Range("M97").Value2 = -75000554
$ws.Range("N97").Value2 = -3527.6

$ws.Range("H116").Value2 = 1576
$ws.Range("I116").Value2 = 864.5
$ws.Range("K116").Value2 = 2593.5
$ws.Range("M116").Value2 = 848.5

$ws.Range("H140").Value2 = 1682.1052
$ws.Range("I140").Value2 = 1682.1052
$ws.Range("K140").Value2 = 5046.3156
$ws.Range("M140").Value2 = 133.6844000000001

$ws.Range("H141").Value2 = 160175.86
$ws.Range("I141").Value2 = 220746.2
$ws.Range("J141").Value2 = 8750
$ws.Range("K141").Value2 = 662238.6000000001
$ws.Range("L141").Value2 = 26250
$ws.Range("M141").Value2 = -657058.6000000001
$ws.Range("N141").Value2 = -36610

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 100002150
$ws.Range("I40").Value2 = 142859090
$ws.Range("J40").Value2 = 2635
$ws.Range("K40").Value2 = 142859090
$ws.Range("L40").Value2 = 2635
$ws.Range("M40").Value2 = -142858954
$ws.Range("N40").Value2 = -2907

$ws.Range("H82").Value2 = 118065
$ws.Range("I82").Value2 = 2749.5
$ws.Range("J82").Value2 = 151012.28
$ws.Range("K82").Value2 = 2749.5
$ws.Range("L82").Value2 = 151012.28
$ws.Range("M82").Value2 = -2388.5
$ws.Range("N82").Value2 = -151734.28

$ws.Range("H85").Value2 = 118065
$ws.Range("I85").Value2 = 2749.5
$ws.Range("J85").Value2 = 151012.28
$ws.Range("K85").Value2 = 2749.5
$ws.Range("L85").Value2 = 151012.28
$ws.Range("M85").Value2 = -1501.5
$ws.Range("N85").Value2 = -153508.28

$ws.Range("H93").Value2 = 20000712
$ws.Range("I93").Value2 = 793.5294
$ws.Range("K93").Value2 = 793.5294
$ws.Range("M93").Value2 = 454.4706

$ws.Range("H122").Value2 = 7144817
$ws.Range("I122").Value2 = 7144817
$ws.Range("K122").Value2 = 21434451
$ws.Range("M122").Value2 = -21432001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value2 = 33672
$ws.Range("J95").Value2 = 33672
$ws.Range("L95").Value2 = 33672
$ws.Range("N95").Value2 = -39164

$ws.Range("H136").Value2 = 4762952
$ws.Range("I136").Value2 = 618.6667
$ws.Range("J136").Value2 = 15153498
$ws.Range("K136").Value2 = 1856.0001
$ws.Range("L136").Value2 = 45460494
$ws.Range("M136").Value2 = 693.9999
$ws.Range("N136").Value2 = -45465594
